$d = $word.ActiveDocument

# Locate the "Frequency Analysis:" heading paragraph, then remove the
# paragraphs that follow it up through the extra blank paragraphs,
# leaving only a single blank paragraph before "Intensity Analysis:".
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Frequency Analysis:") {
        $target = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not find 'Frequency Analysis:' paragraph"
}

# Paragraphs immediately after the heading:
#   target+1 -> "Intuitively the place to start..." text paragraph
#   target+2 -> paragraph containing the inline picture
#   target+3 -> blank paragraph (remove)
#   target+4 -> blank paragraph (remove)
#   target+5 -> blank paragraph (KEEP - becomes the sole spacer paragraph)
$startPara = $d.Paragraphs($target + 1)
$endPara = $d.Paragraphs($target + 4)

$range = $d.Range($startPara.Range.Start, $endPara.Range.End)
$range.Delete()
